$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (row 40) entirely - data now only goes to row 39
$ws.Rows(40).Delete()

# Update B2:C39 with the latest simulation values for the impact parameter
$ws.Cells.Item(2, 2).Value = 432
$ws.Cells.Item(2, 3).Value = 5972.417210512671
$ws.Cells.Item(3, 2).Value = 487
$ws.Cells.Item(3, 3).Value = 6115.046606197327
$ws.Cells.Item(4, 2).Value = 521.5
$ws.Cells.Item(4, 3).Value = 6288.871927893402
$ws.Cells.Item(5, 2).Value = 556
$ws.Cells.Item(5, 3).Value = 6441.067473570095
$ws.Cells.Item(6, 2).Value = 563.5
$ws.Cells.Item(6, 3).Value = 6473.902900933859
$ws.Cells.Item(7, 2).Value = 588.4000000000001
$ws.Cells.Item(7, 3).Value = 6558.985611050061
$ws.Cells.Item(8, 2).Value = 614
$ws.Cells.Item(8, 3).Value = 6659.305107154622
$ws.Cells.Item(9, 2).Value = 639
$ws.Cells.Item(9, 3).Value = 6743.548068764538
$ws.Cells.Item(10, 2).Value = 664.6
$ws.Cells.Item(10, 3).Value = 6849.272070754142
$ws.Cells.Item(11, 2).Value = 689.5
$ws.Cells.Item(11, 3).Value = 6877.455881997951
$ws.Cells.Item(12, 2).Value = 715
$ws.Cells.Item(12, 3).Value = 6945.655292451503
$ws.Cells.Item(13, 2).Value = 740
$ws.Cells.Item(13, 3).Value = 7012.398624499871
$ws.Cells.Item(14, 2).Value = 765
$ws.Cells.Item(14, 3).Value = 7064.864030653095
$ws.Cells.Item(15, 2).Value = 790.5
$ws.Cells.Item(15, 3).Value = 7112.968743412273
$ws.Cells.Item(16, 2).Value = 816
$ws.Cells.Item(16, 3).Value = 7153.876107734214
$ws.Cells.Item(17, 2).Value = 841.3000000000001
$ws.Cells.Item(17, 3).Value = 7212.839593465304
$ws.Cells.Item(18, 2).Value = 866
$ws.Cells.Item(18, 3).Value = 7268.546381995367
$ws.Cells.Item(19, 2).Value = 891.5999999999999
$ws.Cells.Item(19, 3).Value = 7304.881634206284
$ws.Cells.Item(20, 2).Value = 917
$ws.Cells.Item(20, 3).Value = 7319.530341916236
$ws.Cells.Item(21, 2).Value = 942.4
$ws.Cells.Item(21, 3).Value = 7340.992164675764
$ws.Cells.Item(22, 2).Value = 967.3000000000001
$ws.Cells.Item(22, 3).Value = 7378.68323344848
$ws.Cells.Item(23, 2).Value = 992.7
$ws.Cells.Item(23, 3).Value = 7376.375151511461
$ws.Cells.Item(24, 2).Value = 1018
$ws.Cells.Item(24, 3).Value = 7417.355889854926
$ws.Cells.Item(25, 2).Value = 1060
$ws.Cells.Item(25, 3).Value = 7454.962837835755
$ws.Cells.Item(26, 2).Value = 1098
$ws.Cells.Item(26, 3).Value = 7500.377589019557
$ws.Cells.Item(27, 2).Value = 1138
$ws.Cells.Item(27, 3).Value = 7512.81625724276
$ws.Cells.Item(28, 2).Value = 1179
$ws.Cells.Item(28, 3).Value = 7542.425147121889
$ws.Cells.Item(29, 2).Value = 1222
$ws.Cells.Item(29, 3).Value = 7562.234338563289
$ws.Cells.Item(30, 2).Value = 1266
$ws.Cells.Item(30, 3).Value = 7606.5666395152
$ws.Cells.Item(31, 2).Value = 1312
$ws.Cells.Item(31, 3).Value = 7620.387514353233
$ws.Cells.Item(32, 2).Value = 1358
$ws.Cells.Item(32, 3).Value = 7659.663571842446
$ws.Cells.Item(33, 2).Value = 1408
$ws.Cells.Item(33, 3).Value = 7693.906455520149
$ws.Cells.Item(34, 2).Value = 1459
$ws.Cells.Item(34, 3).Value = 7737.625186532094
$ws.Cells.Item(35, 2).Value = 1512
$ws.Cells.Item(35, 3).Value = 7790.108090967696
$ws.Cells.Item(36, 2).Value = 1566
$ws.Cells.Item(36, 3).Value = 7839.396851434865
$ws.Cells.Item(37, 2).Value = 1623
$ws.Cells.Item(37, 3).Value = 7889.153015117722
$ws.Cells.Item(38, 2).Value = 1682
$ws.Cells.Item(38, 3).Value = 7887.856185671338
$ws.Cells.Item(39, 2).Value = 1743
$ws.Cells.Item(39, 3).Value = 7938.074636931969
